$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)
$ws.Range("M16").Insert()
Write-Host "done"
